$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark Status (column C) as DONE for rows 28 and 29
$ws.Range("C28").Value = "DONE"
$ws.Range("C29").Value = "DONE"

# Add new issue row 34 (issue #33)
$ws.Range("A34").Value = 33
$ws.Range("B34").Value = 2
$ws.Range("E34").Value = "icon config"
$ws.Range("H34").Value = "move the config styles for icons into config and icon width"

# Update the selected cell to match the author's final cursor position
$ws.Range("E28").Select() | Out-Null
